# feature: Revising graphs for paper.
# Update numeric values in column B to reflect revised figures used for the paper's graphs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 3.11
$ws.Range("B3").Value  = 95.36499999999999
$ws.Range("B4").Value  = 10655.115
$ws.Range("B10").Value = 0.04434399999999999
$ws.Range("B12").Value = 1.0503985
$ws.Range("B13").Value = 0.05594
$ws.Range("B14").Value = 2.860664
$ws.Range("B15").Value = 1.2350985
$ws.Range("B17").Value = 4.570527687
$ws.Range("B18").Value = 16.55799596
$ws.Range("B19").Value = 0.001536086
$ws.Range("B30").Value = 0.760243
